$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was Melbourne/Brunetti data) becomes the "old" Glen Waverley exposure entry
$ws.Range("A2").Value = "Glen Waverley"
$ws.Range("B2").Value = "Commonwealth Bank, 28-32 Kingsway, Glen Waverley"
$ws.Range("C2").Value = "1:30pm-2:15pm 9/2/2021"
$ws.Range("D2").Value = "Case attended venue"
$ws.Range("E2").Value = "old"

# New row 3: the "new" Glen Waverley exposure entry
$ws.Range("A3").Value = "Glen Waverley"
$ws.Range("B3").Value = "Commonwealth Bank, 28-32 Kingsway, Glen Waverley"
$ws.Range("C3").Value = "1:30pm-2:30pm 9/2/2021"
$ws.Range("D3").Value = "Case attended venue"
$ws.Range("E3").Value = "new"

# New row 4: updated Melbourne / Terminal 4 Melbourne Airport entry
$ws.Range("A4").Value = "Melbourne"
$ws.Range("B4").Value = "Terminal 4, Melbourne Airport"
$ws.Range("C4").Value = "4:45am - 2:00pm  9/2/2021"
$ws.Range("D4").Value = "Case attended venue"
$ws.Range("E4").Value = "new"

# Widen columns A and B to fit the new, longer content (bestFit-style)
$ws.Columns.Item(1).ColumnWidth = 11.29
$ws.Columns.Item(2).ColumnWidth = 43.0

# Match the author's final selection
$ws.Range("B3").Select() | Out-Null
